# Update the simulated-game transition-probability matrix on Sheet1.
# Re-running more simulated games changed the empirical transition
# probabilities in several matrix rows; update the affected cells with
# their newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (state 0)
$ws.Range("B2").Value = 0.1843137254901961
$ws.Range("C2").Value = 0.592156862745098
$ws.Range("J2").Value = 0.007843137254901961
$ws.Range("P2").Value = 0.1333333333333333
$ws.Range("S2").Value = 0.08235294117647059

# Row 3 (state 1)
$ws.Range("B3").Value = 0.00641025641025641
$ws.Range("C3").Value = 0.01923076923076923
$ws.Range("J3").Value = 0.02564102564102564
$ws.Range("P3").Value = 0.6923076923076923
$ws.Range("S3").Value = 0.2564102564102564

# Row 4 (state 2)
$ws.Range("P4").Value = 0.6818181818181818
$ws.Range("S4").Value = 0.3181818181818182

# Row 6 (state 4)
$ws.Range("B6").Value = 0.05785123966942149
$ws.Range("D6").Value = 0.02479338842975207
$ws.Range("F6").Value = 0.07024793388429752
$ws.Range("J6").Value = 0.2479338842975207
$ws.Range("Q6").Value = 0.1570247933884298
$ws.Range("R6").Value = 0.04958677685950413
$ws.Range("S6").Value = 0.3925619834710744

# Row 7 (state 5)
$ws.Range("B7").Value = 0.1149425287356322
$ws.Range("D7").Value = 0.02298850574712644
$ws.Range("E7").Value = 0.005747126436781609
$ws.Range("F7").Value = 0.04022988505747126
$ws.Range("J7").Value = 0.1149425287356322
$ws.Range("O7").Value = 0.01149425287356322
$ws.Range("Q7").Value = 0.1436781609195402
$ws.Range("R7").Value = 0.103448275862069
$ws.Range("S7").Value = 0.4425287356321839

# Row 8 (state 6)
$ws.Range("B8").Value = 0.09562841530054644
$ws.Range("D8").Value = 0.01639344262295082
$ws.Range("E8").Value = 0.00273224043715847
$ws.Range("F8").Value = 0.09836065573770492
$ws.Range("J8").Value = 0.1065573770491803
$ws.Range("O8").Value = 0.01092896174863388
$ws.Range("Q8").Value = 0.2049180327868853
$ws.Range("R8").Value = 0.07103825136612021
$ws.Range("S8").Value = 0.3934426229508197

# Row 9 (state 7)
$ws.Range("B9").Value = 0.07346938775510205
$ws.Range("D9").Value = 0.00816326530612245
$ws.Range("E9").Value = 0.004081632653061225
$ws.Range("F9").Value = 0.0653061224489796
$ws.Range("J9").Value = 0.1020408163265306
$ws.Range("O9").Value = 0.01224489795918367
$ws.Range("Q9").Value = 0.1959183673469388
$ws.Range("R9").Value = 0.08571428571428572
$ws.Range("S9").Value = 0.4530612244897959

# Row 10 (state 8)
$ws.Range("B10").Value = 0.1042893187552565
$ws.Range("D10").Value = 0.02439024390243903
$ws.Range("F10").Value = 0.08494533221194281
$ws.Range("J10").Value = 0.1026072329688814
$ws.Range("O10").Value = 0.009251471825063078
$ws.Range("Q10").Value = 0.2144659377628259
$ws.Range("R10").Value = 0.1093355761143818
$ws.Range("S10").Value = 0.3507148864592094

# Row 11 (state 9)
$ws.Range("G11").Value = 0.1726907630522088
$ws.Range("J11").Value = 0.06024096385542169
$ws.Range("K11").Value = 0.2048192771084337
$ws.Range("L11").Value = 0.5582329317269076
$ws.Range("S11").Value = 0.004016064257028112

# Row 12 (state 10)
$ws.Range("G12").Value = 0.7350993377483444
$ws.Range("J12").Value = 0.1788079470198675
$ws.Range("K12").Value = 0.01324503311258278
$ws.Range("L12").Value = 0.03311258278145696
$ws.Range("S12").Value = 0.03973509933774835

# Row 13 (state 11)
$ws.Range("G13").Value = 0.7058823529411765
$ws.Range("J13").Value = 0.2352941176470588
$ws.Range("S13").Value = 0.05882352941176471

# Row 14 (state 12)
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("S14").Value = 0.3333333333333333

# Row 15 (state 13)
$ws.Range("F15").Value = 0.03296703296703297
$ws.Range("H15").Value = 0.1098901098901099
$ws.Range("I15").Value = 0.1098901098901099
$ws.Range("J15").Value = 0.4230769230769231
$ws.Range("K15").Value = 0.05494505494505494
$ws.Range("O15").Value = 0.07142857142857142
$ws.Range("S15").Value = 0.1978021978021978

# Row 16 (state 14)
$ws.Range("F16").Value = 0.02424242424242424
$ws.Range("H16").Value = 0.1515151515151515
$ws.Range("I16").Value = 0.08484848484848485
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.1696969696969697
$ws.Range("O16").Value = 0.07878787878787878

# Row 17 (state 15)
$ws.Range("F17").Value = 0.009111617312072893
$ws.Range("H17").Value = 0.1708428246013667
$ws.Range("I17").Value = 0.1116173120728929
$ws.Range("J17").Value = 0.4419134396355353
$ws.Range("K17").Value = 0.07289293849658314
$ws.Range("M17").Value = 0.01366742596810934
$ws.Range("N17").Value = 0.002277904328018223
$ws.Range("O17").Value = 0.06605922551252848
$ws.Range("S17").Value = 0.1116173120728929

# Row 18 (state 16)
$ws.Range("F18").Value = 0.01456310679611651
$ws.Range("H18").Value = 0.1699029126213592
$ws.Range("I18").Value = 0.1262135922330097
$ws.Range("J18").Value = 0.412621359223301
$ws.Range("K18").Value = 0.06310679611650485
$ws.Range("M18").Value = 0.02912621359223301
$ws.Range("O18").Value = 0.07281553398058252
$ws.Range("S18").Value = 0.1116504854368932

# Row 19 (state 17)
$ws.Range("F19").Value = 0.02274641954507161
$ws.Range("H19").Value = 0.1735467565290649
$ws.Range("I19").Value = 0.1120471777590564
$ws.Range("J19").Value = 0.3909014321819714
$ws.Range("K19").Value = 0.09098567818028644
$ws.Range("M19").Value = 0.01853411962931761
$ws.Range("N19").Value = 0.001684919966301601
$ws.Range("O19").Value = 0.05728727885425442
$ws.Range("S19").Value = 0.1322662173546756

Write-Host "Updated transition matrix probabilities."
